$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a brand-new paragraph ("30 June 1767 ...") immediately
#    before the existing "2 July 1767" paragraph.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute("2 July 1767", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $savedStart = $anchor.Start
    $anchor.InsertParagraphBefore()

    # The freshly-minted (empty) paragraph now occupies [savedStart, savedStart+1);
    # rebuild a range from the saved offset rather than reusing $anchor (which does
    # not track the shift) or a hard-coded Paragraphs index.
    $newParaRange = $d.Range($savedStart, $savedStart + 1)

    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b/><w:color w:val="000000"/></w:rPr><w:t>30 June 1767</w:t></w:r>' + `
        '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">  This is the approximate date of the first performance of </w:t></w:r>' + `
        '<w:r><w:rPr><w:i/><w:color w:val="000000"/></w:rPr><w:t>Lucio Papirio Dittatore</w:t></w:r>' + `
        '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>, a dramma per musica by Giovanni Paisiello (27) to words of Zeno, in Teatro San Carlo, Naples.</w:t></w:r>' + `
        '</w:p>'

    $newParaRange.InsertXML($newParaXml)
}

# ---------------------------------------------------------------------
# 2) Bump the copyright year: 2004-2015 -> 2004-2016
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Paul Scharfenberger 2004-2015", $true, $false, $false, $false, $false, $true, 1, $false, "Paul Scharfenberger 2004-2016", 2)

# ---------------------------------------------------------------------
# 3) Update the revision date paragraph: "2 December 2015" -> "2 June 2016"
#    Keep the leading "2 " as its own run and fold the rest into a single
#    "June 2016" run, matching the target run layout exactly.
# ---------------------------------------------------------------------
$dateAnchor = $d.Content
$dateAnchor.Find.ClearFormatting()
$dateFound = $dateAnchor.Find.Execute("2 December 2015", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($dateFound) {
    $dateAnchor.Expand(4)
    $dateXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:r><w:t xml:space="preserve">2 </w:t></w:r>' + `
        '<w:r><w:t>June 2016</w:t></w:r>' + `
        '</w:p>'
    $dateAnchor.InsertXML($dateXml)
}
